$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: columns E..O get new/shifted headers
$ws.Range("E1").Value = "INTEREST.RATE"
$ws.Range("F1").Value = "INTEND.DATE"
$ws.Range("G1").Value = "CUST.REMARKS:1"
$ws.Range("H1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("I1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("J1").Value = "PRIN.LIQ.ACCT"
$ws.Range("K1").Value = "INT.LIQ.ACCT"
$ws.Range("L1").Value = "CHRG.LIQ.ACCT"
$ws.Range("M1").Value = "AUTO.ROLLOVER"
$ws.Range("N1").Value = "FINAL.MATURITY"
$ws.Range("O1").Value = "EXP.DATE"

# Move the AUTO.ROLLOVER data value from E2 to M2
$ws.Range("E2").ClearContents()
$ws.Range("M2").Value = 2

# Update the active selection to match the saved view state
$ws.Range("J9").Select()
